$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 7400
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 10600
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 10600
$ws.Range("M21").Value = -532
$ws.Range("N21").Value = -11536
$ws.Range("H23").Value = 7400
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 10600
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 10600
$ws.Range("M23").Value = -766
$ws.Range("N23").Value = -11068
$ws.Range("H32").Value = 2500
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2500
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2500
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3152
$ws.Range("H80").Value = 682.8889
$ws.Range("I80").Value = 278.8
$ws.Range("J80").Value = 1188
$ws.Range("K80").Value = 836.4000000000001
$ws.Range("L80").Value = 3564
$ws.Range("M80").Value = 161.5999999999999
$ws.Range("N80").Value = -5560
$ws.Range("H83").Value = 682.8889
$ws.Range("I83").Value = 278.8
$ws.Range("J83").Value = 1188
$ws.Range("K83").Value = 2509.2
$ws.Range("L83").Value = 10692
$ws.Range("M83").Value = 2482.8
$ws.Range("N83").Value = -20676
$ws.Range("H132").Value = 2870.976
$ws.Range("I132").Value = 2828.8857
$ws.Range("J132").Value = 3081.4285
$ws.Range("K132").Value = 8486.6571
$ws.Range("L132").Value = 9244.2855
$ws.Range("M132").Value = -5956.6571
$ws.Range("N132").Value = -14304.2855
$ws.Range("H137").Value = 2165
$ws.Range("I137").Value = 1723
$ws.Range("K137").Value = 5169
$ws.Range("M137").Value = -2619
$ws.Range("H138").Value = 3200
$ws.Range("I138").Value = 1498.5
$ws.Range("K138").Value = 4495.5
$ws.Range("M138").Value = 644.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 846.25
$ws.Range("I5").Value = 961.6667
$ws.Range("K5").Value = 961.6667
$ws.Range("M5").Value = -849.6667
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H63").Value = 2714.7
$ws.Range("I63").Value = 2030.2
$ws.Range("J63").Value = 3399.2
$ws.Range("K63").Value = 2030.2
$ws.Range("L63").Value = 3399.2
$ws.Range("M63").Value = -1344.2
$ws.Range("N63").Value = -4771.2
$ws.Range("H66").Value = 2714.7
$ws.Range("I66").Value = 2030.2
$ws.Range("J66").Value = 3399.2
$ws.Range("K66").Value = 10151
$ws.Range("L66").Value = 16996
$ws.Range("M66").Value = -6719
$ws.Range("N66").Value = -23860
$ws.Range("H74").Value = 1421.6
$ws.Range("I74").Value = 1209.6154
$ws.Range("K74").Value = 1209.6154
$ws.Range("M74").Value = -335.6153999999999
$ws.Range("H77").Value = 1421.6
$ws.Range("I77").Value = 1209.6154
$ws.Range("K77").Value = 6048.076999999999
$ws.Range("M77").Value = -1680.076999999999
$ws.Range("H88").Value = 9398.200000000001
$ws.Range("J88").Value = 9398.200000000001
$ws.Range("L88").Value = 9398.200000000001
$ws.Range("N88").Value = -10210.2
$ws.Range("H91").Value = 9398.200000000001
$ws.Range("J91").Value = 9398.200000000001
$ws.Range("L91").Value = 9398.200000000001
$ws.Range("N91").Value = -12206.2

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 846.25
$ws.Range("I4").Value = 961.6667
$ws.Range("K4").Value = 961.6667
$ws.Range("M4").Value = -846.6667

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5722.222
$ws.Range("I58").Value = 8762.5
$ws.Range("J58").Value = 3290
$ws.Range("K58").Value = 8762.5
$ws.Range("L58").Value = 3290
$ws.Range("M58").Value = -8559.5
$ws.Range("N58").Value = -3696
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H80").Value = 10800
$ws.Range("J80").Value = 10800
$ws.Range("L80").Value = 10800
$ws.Range("N80").Value = -13046
$ws.Range("H83").Value = 10800
$ws.Range("J83").Value = 10800
$ws.Range("L83").Value = 32400
$ws.Range("N83").Value = -43632
$ws.Range("H88").Value = 11345.077
$ws.Range("J88").Value = 11345.077
$ws.Range("L88").Value = 11345.077
$ws.Range("N88").Value = -12157.077
$ws.Range("H91").Value = 11345.077
$ws.Range("J91").Value = 11345.077
$ws.Range("L91").Value = 11345.077
$ws.Range("N91").Value = -14153.077
$ws.Range("H103").Value = 8731
$ws.Range("I103").Value = 8731
$ws.Range("K103").Value = 8731
$ws.Range("M103").Value = -7559
$ws.Range("H122").Value = 4016.8572
$ws.Range("I122").Value = 5417.125
$ws.Range("J122").Value = 2149.8333
$ws.Range("K122").Value = 16251.375
$ws.Range("L122").Value = 6449.499899999999
$ws.Range("M122").Value = -13801.375
$ws.Range("N122").Value = -11349.4999
$ws.Range("H132").Value = 3279
$ws.Range("I132").Value = 3279
$ws.Range("K132").Value = 9837
$ws.Range("M132").Value = -7307
$ws.Range("H134").Value = 1966.35
$ws.Range("I134").Value = 2176.4375
$ws.Range("J134").Value = 1126
$ws.Range("K134").Value = 6529.3125
$ws.Range("L134").Value = 3378
$ws.Range("M134").Value = -3994.3125
$ws.Range("N134").Value = -8448
$ws.Range("H136").Value = 5722.222
$ws.Range("I136").Value = 8762.5
$ws.Range("J136").Value = 3290
$ws.Range("K136").Value = 26287.5
$ws.Range("L136").Value = 9870
$ws.Range("M136").Value = -23737.5
$ws.Range("N136").Value = -14970

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 5780
$ws.Range("J22").Value = 5780
$ws.Range("L22").Value = 17340
$ws.Range("N22").Value = -17678
$ws.Range("H27").Value = 5780
$ws.Range("J27").Value = 5780
$ws.Range("L27").Value = 17340
$ws.Range("N27").Value = -17544
$ws.Range("H29").Value = 73.84614999999999
$ws.Range("I29").Value = 190
$ws.Range("J29").Value = 64.166664
$ws.Range("K29").Value = 570
$ws.Range("L29").Value = 192.499992
$ws.Range("M29").Value = -293
$ws.Range("N29").Value = -746.499992
$ws.Range("H56").Value = 17854.615
$ws.Range("I56").Value = 17854.615
$ws.Range("K56").Value = 17854.615
$ws.Range("M56").Value = -17324.615
$ws.Range("H119").Value = 608.75
$ws.Range("I119").Value = 608.75
$ws.Range("K119").Value = 1826.25
$ws.Range("M119").Value = 3011.75
$ws.Range("H122").Value = 1852
$ws.Range("I122").Value = 1796.4
$ws.Range("J122").Value = 1921.5
$ws.Range("K122").Value = 16167.6
$ws.Range("L122").Value = 17293.5
$ws.Range("M122").Value = -13717.6
$ws.Range("N122").Value = -22193.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 24999.334
$ws.Range("J44").Value = 24999.334
$ws.Range("L44").Value = 24999.334
$ws.Range("N44").Value = -26191.334
$ws.Range("H122").Value = 39448.21
$ws.Range("I122").Value = 38094.812
$ws.Range("K122").Value = 114284.436
$ws.Range("M122").Value = -111834.436

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 9999
$ws.Range("J5").Value = 9999
$ws.Range("L5").Value = 9999
$ws.Range("N5").Value = -10225
$ws.Range("H7").Value = 1858
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H126").Value = 1858
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H127").Value = 57500
$ws.Range("J127").Value = 57500
$ws.Range("L127").Value = 57500
$ws.Range("N127").Value = -67420
$ws.Range("H132").Value = 14573.966
$ws.Range("J132").Value = 10393.833
$ws.Range("L132").Value = 31181.499
$ws.Range("N132").Value = -36241.499

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 35000
$ws.Range("J54").Value = 35000
$ws.Range("L54").Value = 35000
$ws.Range("N54").Value = -36040
$ws.Range("H132").Value = 499.5
$ws.Range("I132").Value = 499.5
$ws.Range("K132").Value = 1498.5
$ws.Range("M132").Value = 1031.5
$ws.Range("H136").Value = 3923.0908
$ws.Range("I136").Value = 3923.0908
$ws.Range("K136").Value = 11769.2724
$ws.Range("M136").Value = -9219.2724
